$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 106.416664
$ws.Range("I4").Value = 98.625
$ws.Range("J4").Value = 122
$ws.Range("K4").Value = 98.625
$ws.Range("L4").Value = 122
$ws.Range("M4").Value = 15.375
$ws.Range("N4").Value = -350
$ws.Range("H8").Value = 3000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 9000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -9278
$ws.Range("H18").Value = 3217
$ws.Range("I18").Value = 1860.4
$ws.Range("K18").Value = 1860.4
$ws.Range("M18").Value = -1576.4
$ws.Range("H80").Value = 423
$ws.Range("I80").Value = 495
$ws.Range("K80").Value = 1485
$ws.Range("M80").Value = -487
$ws.Range("H83").Value = 423
$ws.Range("I83").Value = 495
$ws.Range("K83").Value = 4455
$ws.Range("M83").Value = 537
$ws.Range("H92").Value = 1172
$ws.Range("I92").Value = 1225.85
$ws.Range("J92").Value = 95
$ws.Range("K92").Value = 1225.85
$ws.Range("L92").Value = 95
$ws.Range("M92").Value = 22.15000000000009
$ws.Range("N92").Value = -2591
$ws.Range("H98").Value = 1446.9
$ws.Range("I98").Value = 829.1739
$ws.Range("J98").Value = 3476.5715
$ws.Range("K98").Value = 829.1739
$ws.Range("L98").Value = 3476.5715
$ws.Range("M98").Value = 668.8261
$ws.Range("N98").Value = -6472.5715
$ws.Range("H99").Value = 408.625
$ws.Range("I99").Value = 432
$ws.Range("K99").Value = 1296
$ws.Range("M99").Value = 202
$ws.Range("H101").Value = 1236
$ws.Range("I101").Value = 165.5
$ws.Range("K101").Value = 496.5
$ws.Range("M101").Value = 1125.5
$ws.Range("H122").Value = 1446.9
$ws.Range("I122").Value = 829.1739
$ws.Range("J122").Value = 3476.5715
$ws.Range("K122").Value = 2487.5217
$ws.Range("L122").Value = 10429.7145
$ws.Range("M122").Value = -37.52170000000024
$ws.Range("N122").Value = -15329.7145
$ws.Range("H132").Value = 1442.5333
$ws.Range("I132").Value = 1447.738
$ws.Range("J132").Value = 1369.6666
$ws.Range("K132").Value = 4343.214
$ws.Range("L132").Value = 4108.9998
$ws.Range("M132").Value = -1813.214
$ws.Range("N132").Value = -9168.9998
$ws.Range("H135").Value = 1305.6552
$ws.Range("I135").Value = 894.9583
$ws.Range("K135").Value = 8054.6247
$ws.Range("M135").Value = -5519.6247
$ws.Range("H137").Value = 1546.25
$ws.Range("I137").Value = 1409.762
$ws.Range("J137").Value = 1737.3334
$ws.Range("K137").Value = 4229.286
$ws.Range("L137").Value = 5212.0002
$ws.Range("M137").Value = -1679.286
$ws.Range("N137").Value = -10312.0002
$ws.Range("H138").Value = 6852249
$ws.Range("I138").Value = 1636.4546
$ws.Range("K138").Value = 4909.3638
$ws.Range("M138").Value = 230.6361999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7426.0806
$ws.Range("I32").Value = 3182.325
$ws.Range("J32").Value = 25294.525
$ws.Range("K32").Value = 3182.325
$ws.Range("L32").Value = 25294.525
$ws.Range("M32").Value = -2895.325
$ws.Range("N32").Value = -25868.525
$ws.Range("H74").Value = 7157.1816
$ws.Range("I74").Value = 1399.7333
$ws.Range("K74").Value = 1399.7333
$ws.Range("M74").Value = -525.7333000000001
$ws.Range("H77").Value = 7157.1816
$ws.Range("I77").Value = 1399.7333
$ws.Range("K77").Value = 6998.6665
$ws.Range("M77").Value = -2630.6665
$ws.Range("H132").Value = 2994
$ws.Range("I132").Value = 2747.0454
$ws.Range("J132").Value = 4352.25
$ws.Range("K132").Value = 8241.136200000001
$ws.Range("L132").Value = 13056.75
$ws.Range("M132").Value = -5711.136200000001
$ws.Range("N132").Value = -18116.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6395.1763
$ws.Range("I86").Value = 6159.154
$ws.Range("J86").Value = 7162.25
$ws.Range("K86").Value = 6159.154
$ws.Range("L86").Value = 7162.25
$ws.Range("M86").Value = -5036.154
$ws.Range("N86").Value = -9408.25
$ws.Range("H89").Value = 6395.1763
$ws.Range("I89").Value = 6159.154
$ws.Range("J89").Value = 7162.25
$ws.Range("K89").Value = 30795.77
$ws.Range("L89").Value = 35811.25
$ws.Range("M89").Value = -25179.77
$ws.Range("N89").Value = -47043.25
$ws.Range("H99").Value = 230286.44
$ws.Range("I99").Value = 339999.34
$ws.Range("J99").Value = 175430
$ws.Range("K99").Value = 339999.34
$ws.Range("L99").Value = 175430
$ws.Range("M99").Value = -338501.34
$ws.Range("N99").Value = -178426
$ws.Range("H134").Value = 2115.8193
$ws.Range("I134").Value = 2005.3896
$ws.Range("K134").Value = 6016.168799999999
$ws.Range("M134").Value = -3481.168799999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1749.9
$ws.Range("I3").Value = 1350
$ws.Range("J3").Value = 2016.5
$ws.Range("K3").Value = 1350
$ws.Range("L3").Value = 2016.5
$ws.Range("M3").Value = -1237
$ws.Range("N3").Value = -2242.5
$ws.Range("H7").Value = 405.69446
$ws.Range("I7").Value = 496.5
$ws.Range("K7").Value = 496.5
$ws.Range("M7").Value = -383.5
$ws.Range("H31").Value = 74032.86
$ws.Range("I31").Value = 92799.17999999999
$ws.Range("J31").Value = 5223
$ws.Range("K31").Value = 92799.17999999999
$ws.Range("L31").Value = 5223
$ws.Range("M31").Value = -92504.17999999999
$ws.Range("N31").Value = -5813
$ws.Range("H34").Value = 74032.86
$ws.Range("I34").Value = 92799.17999999999
$ws.Range("J34").Value = 5223
$ws.Range("K34").Value = 92799.17999999999
$ws.Range("L34").Value = 5223
$ws.Range("M34").Value = -92597.17999999999
$ws.Range("N34").Value = -5627
$ws.Range("H37").Value = 60000
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H58").Value = 2899.2942
$ws.Range("I58").Value = 1700.0834
$ws.Range("K58").Value = 1700.0834
$ws.Range("M58").Value = -1497.0834
$ws.Range("H86").Value = 7772.0713
$ws.Range("I86").Value = 7662.3076
$ws.Range("K86").Value = 7662.3076
$ws.Range("M86").Value = -6539.3076
$ws.Range("H89").Value = 7772.0713
$ws.Range("I89").Value = 7662.3076
$ws.Range("K89").Value = 38311.538
$ws.Range("M89").Value = -32695.538
$ws.Range("H132").Value = 3995.5715
$ws.Range("I132").Value = 3962.25
$ws.Range("K132").Value = 11886.75
$ws.Range("M132").Value = -9356.75
$ws.Range("H134").Value = 8733.093999999999
$ws.Range("I134").Value = 5342.7915
$ws.Range("J134").Value = 18904
$ws.Range("K134").Value = 16028.3745
$ws.Range("L134").Value = 56712
$ws.Range("M134").Value = -13493.3745
$ws.Range("N134").Value = -61782
$ws.Range("H136").Value = 2899.2942
$ws.Range("I136").Value = 1700.0834
$ws.Range("K136").Value = 5100.2502
$ws.Range("M136").Value = -2550.2502

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 15.964286
$ws.Range("I2").Value = 20
$ws.Range("K2").Value = 120
$ws.Range("M2").Value = -7
$ws.Range("H51").Value = 967
$ws.Range("I51").Value = 967
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 2901
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -2441
$ws.Range("N51").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 60000
$ws.Range("J38").Value = 60000
$ws.Range("L38").Value = 60000
$ws.Range("N38").Value = -60926
$ws.Range("H102").Value = 1000000000
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H126").Value = 21425.75
$ws.Range("I126").Value = 27167.777
$ws.Range("K126").Value = 81503.33099999999
$ws.Range("M126").Value = -79033.33099999999
$ws.Range("H132").Value = 2490.25
$ws.Range("I132").Value = 1908.3846
$ws.Range("J132").Value = 5011.6665
$ws.Range("K132").Value = 5725.1538
$ws.Range("L132").Value = 15034.9995
$ws.Range("M132").Value = -3195.1538
$ws.Range("N132").Value = -20094.9995
$ws.Range("H141").Value = 90416.664
$ws.Range("J141").Value = 90416.664
$ws.Range("L141").Value = 90416.664
$ws.Range("N141").Value = -100776.664

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 802.3333
$ws.Range("I46").Value = 575.3333
$ws.Range("J46").Value = 1483.3334
$ws.Range("K46").Value = 575.3333
$ws.Range("L46").Value = 1483.3334
$ws.Range("M46").Value = -387.3333
$ws.Range("N46").Value = -1859.3334
$ws.Range("H122").Value = 4996.933
$ws.Range("I122").Value = 4917.1113
$ws.Range("K122").Value = 14751.3339
$ws.Range("M122").Value = -12301.3339
$ws.Range("H136").Value = 4857.5713
$ws.Range("I136").Value = 4167.1665
$ws.Range("K136").Value = 12501.4995
$ws.Range("M136").Value = -9951.499500000002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4722.1113
$ws.Range("J4").Value = 4937.375
$ws.Range("L4").Value = 4937.375
$ws.Range("N4").Value = -5163.375
$ws.Range("H126").Value = 3999.8
$ws.Range("J126").Value = 4000
$ws.Range("L126").Value = 12000
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 3154.0435
$ws.Range("I132").Value = 2776.4211
$ws.Range("K132").Value = 8329.263300000001
$ws.Range("M132").Value = -5799.263300000001
$ws.Range("H136").Value = 34928.066
$ws.Range("I136").Value = 51062.7
$ws.Range("K136").Value = 153188.1
$ws.Range("M136").Value = -150638.1
$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360
$ws.Range("H141").Value = 64454.453
$ws.Range("J141").Value = 64454.453
$ws.Range("L141").Value = 64454.453
$ws.Range("N141").Value = -74814.45300000001
